$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testcases")

# The "infor_TC*_...2" block (rows 8-10) is duplicated twice below it to add
# a "3" suffixed block (rows 11-13) and a "4" suffixed block (rows 14-16),
# giving extra test data for parallel execution.
$ws.Range("A8:M10").Copy()
$ws.Paste($ws.Range("A11:M13"))
$ws.Range("A8:M10").Copy()
$ws.Paste($ws.Range("A14:M16"))

# Fix up the test-case name column for the newly pasted blocks.
$ws.Range("A11").Value = "infor_TC1_Login3"
$ws.Range("A12").Value = "infor_TC2_Registration3"
$ws.Range("A13").Value = "infor_TC3_LoginViaConfigFile3"

$ws.Range("A14").Value = "infor_TC1_Login4"
$ws.Range("A15").Value = "infor_TC2_Registration4"
$ws.Range("A16").Value = "infor_TC3_LoginViaConfigFile4"

# Restore the leading-zero phone number text (quote-prefixed) in the
# newly duplicated registration rows.
$ws.Range("G12").Value = "'09661401029"
$ws.Range("G15").Value = "'09661401029"

# Add the hyperlinks for the email cells in the newly created registration
# rows, then restore the Hyperlink cell style (Add() re-styles the cell
# with its own font-only style, so reapply the named "Hyperlink" style).
$ws.Hyperlinks.Add($ws.Range("F12"), "mailto:jazx.zn@gmail.com") | Out-Null
$ws.Range("F12").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F15"), "mailto:jazx.zn@gmail.com") | Out-Null
$ws.Range("F15").Style = "Hyperlink"

$ws.Range("A17").Select()
